$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 775.2632
$ws.Range("J19").Value = 818.2727
$ws.Range("L19").Value = 818.2727
$ws.Range("N19").Value = -1168.2727

$ws.Range("H129").Value = 980.60315
$ws.Range("J129").Value = 1024.7069
$ws.Range("L129").Value = 3074.120699999999
$ws.Range("N129").Value = -13074.1207

$ws.Range("H138").Value = 9806574
$ws.Range("I138").Value = 3070.0557
$ws.Range("J138").Value = 15153939
$ws.Range("K138").Value = 9210.167099999999
$ws.Range("L138").Value = 45461817
$ws.Range("M138").Value = -4070.167099999999
$ws.Range("N138").Value = -45472097

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5199.5557
$ws.Range("I2").Value = 6636
$ws.Range("J2").Value = 1095.4286
$ws.Range("K2").Value = 6636
$ws.Range("L2").Value = 1095.4286
$ws.Range("M2").Value = -6523
$ws.Range("N2").Value = -1321.4286

$ws.Range("H32").Value = 2211.6858
$ws.Range("I32").Value = 1764.9788
$ws.Range("J32").Value = 3124.5217
$ws.Range("K32").Value = 1764.9788
$ws.Range("L32").Value = 3124.5217
$ws.Range("M32").Value = -1477.9788
$ws.Range("N32").Value = -3698.5217

$ws.Range("H74").Value = 696.9737
$ws.Range("I74").Value = 641.7222
$ws.Range("J74").Value = 746.7
$ws.Range("K74").Value = 641.7222
$ws.Range("L74").Value = 746.7
$ws.Range("M74").Value = 232.2778
$ws.Range("N74").Value = -2494.7

$ws.Range("H77").Value = 696.9737
$ws.Range("I77").Value = 641.7222
$ws.Range("J77").Value = 746.7
$ws.Range("K77").Value = 3208.611
$ws.Range("L77").Value = 3733.5
$ws.Range("M77").Value = 1159.389
$ws.Range("N77").Value = -12469.5

$ws.Range("H110").Value = 1158.1
$ws.Range("I110").Value = 893.3333
$ws.Range("J110").Value = 1555.25
$ws.Range("K110").Value = 893.3333
$ws.Range("L110").Value = 1555.25
$ws.Range("M110").Value = 1151.6667
$ws.Range("N110").Value = -5645.25

$ws.Range("H116").Value = 5199.5557
$ws.Range("I116").Value = 6636
$ws.Range("J116").Value = 1095.4286
$ws.Range("K116").Value = 6636
$ws.Range("L116").Value = 1095.4286
$ws.Range("M116").Value = -4342
$ws.Range("N116").Value = -5683.4286

$ws.Range("H139").Value = 58812.57
$ws.Range("J139").Value = 58812.57
$ws.Range("L139").Value = 58812.57
$ws.Range("N139").Value = -69092.57000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5199.5557
$ws.Range("I3").Value = 6636
$ws.Range("J3").Value = 1095.4286
$ws.Range("K3").Value = 6636
$ws.Range("L3").Value = 1095.4286
$ws.Range("M3").Value = -6522
$ws.Range("N3").Value = -1323.4286

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 10574.917
$ws.Range("J39").Value = 10947.739
$ws.Range("L39").Value = 32843.217
$ws.Range("N39").Value = -33431.217

$ws.Range("H68").Value = 1746.6305
$ws.Range("I68").Value = 1618.7812
$ws.Range("J68").Value = 2038.8572
$ws.Range("K68").Value = 4856.3436
$ws.Range("L68").Value = 6116.571599999999
$ws.Range("M68").Value = -4045.3436
$ws.Range("N68").Value = -7738.571599999999

$ws.Range("H71").Value = 1746.6305
$ws.Range("I71").Value = 1618.7812
$ws.Range("J71").Value = 2038.8572
$ws.Range("K71").Value = 14569.0308
$ws.Range("L71").Value = 18349.7148
$ws.Range("M71").Value = -10513.0308
$ws.Range("N71").Value = -26461.7148

$ws.Range("H98").Value = 325.66666
$ws.Range("I98").Value = 350
$ws.Range("J98").Value = 204
$ws.Range("K98").Value = 1050
$ws.Range("L98").Value = 612
$ws.Range("M98").Value = 448
$ws.Range("N98").Value = -3608

$ws.Range("H106").Value = 3698.375
$ws.Range("J106").Value = 3869.5715
$ws.Range("L106").Value = 11608.7145
$ws.Range("N106").Value = -13500.7145

$ws.Range("H109").Value = 901
$ws.Range("I109").Value = 901
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 2703
$ws.Range("L109").ClearContents()
$ws.Range("M109").Value = -1663
$ws.Range("N109").Value = 0

$ws.Range("H112").Value = 125002290
$ws.Range("I112").Value = 925.6667
$ws.Range("J112").Value = 200003100
$ws.Range("K112").Value = 2777.0001
$ws.Range("L112").Value = 600009300
$ws.Range("M112").Value = -1669.0001
$ws.Range("N112").Value = -600011516

$ws.Range("H113").Value = 396
$ws.Range("I113").Value = 413.5
$ws.Range("J113").Value = 394.54166
$ws.Range("K113").Value = 1240.5
$ws.Range("L113").Value = 1183.62498
$ws.Range("M113").Value = 929.5
$ws.Range("N113").Value = -5523.624980000001

$ws.Range("H115").Value = 2084.1428
$ws.Range("J115").Value = 2747.5
$ws.Range("L115").Value = 8242.5
$ws.Range("N115").Value = -10592.5

$ws.Range("H118").Value = 1322.8
$ws.Range("I118").Value = 853.5
$ws.Range("J118").Value = 3200
$ws.Range("K118").Value = 2560.5
$ws.Range("L118").Value = 9600
$ws.Range("M118").Value = -1317.5
$ws.Range("N118").Value = -12086

$ws.Range("H121").Value = 0
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = 0
$ws.Range("L121").ClearContents()
$ws.Range("M121").ClearContents()
$ws.Range("N121").Value = 0

$ws.Range("H131").Value = 2743.3428
$ws.Range("J131").Value = 3045.7097
$ws.Range("L131").Value = 9137.1291
$ws.Range("N131").Value = -19217.1291

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H104").Value = 15671
$ws.Range("J104").Value = 15671
$ws.Range("L104").Value = 15671
$ws.Range("N104").Value = -22659

$ws.Range("H122").Value = 3710370.2
$ws.Range("I122").Value = 3710370.2
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 11131110.6
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -11128660.6

$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").ClearContents()
$ws.Range("N128").Value = 0

$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("L130").ClearContents()
$ws.Range("M130").ClearContents()
$ws.Range("N130").Value = 0

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3416.6667
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()

$ws.Range("H122").Value = 3210.5862
$ws.Range("I122").Value = 2075.875
$ws.Range("J122").Value = 3642.8572
$ws.Range("K122").Value = 6227.625
$ws.Range("L122").Value = 10928.5716
$ws.Range("M122").Value = -3777.625
$ws.Range("N122").Value = -15828.5716

$ws.Range("H126").Value = 3416.6667
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H57").Value = 51000
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()

$ws.Range("H118").Value = 30000
$ws.Range("J118").Value = 30000
$ws.Range("L118").Value = 30000
$ws.Range("N118").Value = -33314

$ws.Range("H121").Value = 26592
$ws.Range("J121").Value = 26592
$ws.Range("L121").Value = 26592
$ws.Range("N121").Value = -30086

$ws.Range("H123").Value = 35429
$ws.Range("J123").Value = 35429
$ws.Range("L123").Value = 35429
$ws.Range("N123").Value = -45229

$ws.Range("H125").Value = 34000
$ws.Range("J125").Value = 34000
$ws.Range("L125").Value = 34000
$ws.Range("N125").Value = -43840

$ws.Range("H128").Value = 59800
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 59800
